# Auto-generated edit script applying the Hades_Profits.xlsx diff
# (workbook = Sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 404.2353
$ws.Range("I33").Value = 384.76923
$ws.Range("K33").Value = 384.76923
$ws.Range("M33").Value = -155.76923

# Row 64
$ws.Range("H64").Value = 4093.0386
$ws.Range("I64").Value = 3737.8462
$ws.Range("J64").Value = 4448.231
$ws.Range("K64").Value = 3737.8462
$ws.Range("L64").Value = 4448.231
$ws.Range("M64").Value = -3489.8462
$ws.Range("N64").Value = -4944.231

# Row 67
$ws.Range("H67").Value = 4093.0386
$ws.Range("I67").Value = 3737.8462
$ws.Range("J67").Value = 4448.231
$ws.Range("K67").Value = 3737.8462
$ws.Range("L67").Value = 4448.231
$ws.Range("M67").Value = -2879.8462
$ws.Range("N67").Value = -6164.231

# Row 70
$ws.Range("H70").Value = 3510.4
$ws.Range("I70").Value = 1475.25
$ws.Range("J70").Value = 4867.1665
$ws.Range("K70").Value = 4425.75
$ws.Range("L70").Value = 14601.4995
$ws.Range("M70").Value = -4155.75
$ws.Range("N70").Value = -15141.4995

# Row 73
$ws.Range("H73").Value = 3510.4
$ws.Range("I73").Value = 1475.25
$ws.Range("J73").Value = 4867.1665
$ws.Range("K73").Value = 4425.75
$ws.Range("L73").Value = 14601.4995
$ws.Range("M73").Value = -3489.75
$ws.Range("N73").Value = -16473.4995

# Row 74
$ws.Range("H74").Value = 3600.6924
$ws.Range("I74").Value = 3527.3635
$ws.Range("J74").Value = 4004
$ws.Range("K74").Value = 3527.3635
$ws.Range("L74").Value = 4004
$ws.Range("M74").Value = -2591.3635
$ws.Range("N74").Value = -5876

# Row 77
$ws.Range("H77").Value = 3600.6924
$ws.Range("I77").Value = 3527.3635
$ws.Range("J77").Value = 4004
$ws.Range("K77").Value = 17636.8175
$ws.Range("L77").Value = 20020
$ws.Range("M77").Value = -12956.8175
$ws.Range("N77").Value = -29380

# Row 134
$ws.Range("H134").Value = 60000
$ws.Range("J134").Value = 60000
$ws.Range("L134").Value = 60000
$ws.Range("N134").Value = -70140

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7012.04
$ws.Range("I32").Value = 4163.782
$ws.Range("J32").Value = 17110.408
$ws.Range("K32").Value = 4163.782
$ws.Range("L32").Value = 17110.408
$ws.Range("M32").Value = -3876.782
$ws.Range("N32").Value = -17684.408

# Row 45
$ws.Range("H45").Value = 2170.6667
$ws.Range("I45").Value = 2406
$ws.Range("J45").Value = 1700
$ws.Range("K45").Value = 2406
$ws.Range("L45").Value = 1700
$ws.Range("M45").Value = -2029
$ws.Range("N45").Value = -2454

# Row 61
$ws.Range("H61").Value = 12220932
$ws.Range("I61").Value = 13172268
$ws.Range("J61").Value = 170669
$ws.Range("K61").Value = 13172268
$ws.Range("L61").Value = 170669
$ws.Range("M61").Value = -13172056
$ws.Range("N61").Value = -171093

# Row 132
$ws.Range("H132").Value = 57046.92
$ws.Range("I132").Value = 45358.043
$ws.Range("J132").Value = 74969.87
$ws.Range("K132").Value = 136074.129
$ws.Range("L132").Value = 224909.61
$ws.Range("M132").Value = -133544.129
$ws.Range("N132").Value = -229969.61

# Row 136
$ws.Range("H136").Value = 12220932
$ws.Range("I136").Value = 13172268
$ws.Range("J136").Value = 170669
$ws.Range("K136").Value = 39516804
$ws.Range("L136").Value = 512007
$ws.Range("M136").Value = -39514254
$ws.Range("N136").Value = -517107

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2567.8462
$ws.Range("I134").Value = 2567.8462
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7703.5386
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -5168.5386
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6832.242
$ws.Range("I31").Value = 2674.5454
$ws.Range("J31").Value = 8911.091
$ws.Range("K31").Value = 2674.5454
$ws.Range("L31").Value = 8911.091
$ws.Range("M31").Value = -2379.5454
$ws.Range("N31").Value = -9501.091

# Row 34
$ws.Range("H34").Value = 6832.242
$ws.Range("I34").Value = 2674.5454
$ws.Range("J34").Value = 8911.091
$ws.Range("K34").Value = 2674.5454
$ws.Range("L34").Value = 8911.091
$ws.Range("M34").Value = -2472.5454
$ws.Range("N34").Value = -9315.091

# Row 134
$ws.Range("H134").Value = 23725.584
$ws.Range("I134").Value = 1184.1025
$ws.Range("J134").Value = 121405.336
$ws.Range("K134").Value = 3552.3075
$ws.Range("L134").Value = 364216.008
$ws.Range("M134").Value = -1017.3075
$ws.Range("N134").Value = -369286.008

$ws = $wb.Worksheets.Item("CUL")
# Row 75
$ws.Range("H75").Value = 2343.125
$ws.Range("I75").Value = 1732.75
$ws.Range("J75").Value = 2953.5
$ws.Range("K75").Value = 5198.25
$ws.Range("L75").Value = 8860.5
$ws.Range("M75").Value = -4200.25
$ws.Range("N75").Value = -10856.5

# Row 78
$ws.Range("H78").Value = 2343.125
$ws.Range("I78").Value = 1732.75
$ws.Range("J78").Value = 2953.5
$ws.Range("K78").Value = 15594.75
$ws.Range("L78").Value = 26581.5
$ws.Range("M78").Value = -10602.75
$ws.Range("N78").Value = -36565.5

# Row 87
$ws.Range("H87").Value = 29449.875
$ws.Range("I87").Value = 11700
$ws.Range("J87").Value = 32999.85
$ws.Range("K87").Value = 35100
$ws.Range("L87").Value = 98999.54999999999
$ws.Range("M87").Value = -33852
$ws.Range("N87").Value = -101495.55

# Row 90
$ws.Range("H90").Value = 29449.875
$ws.Range("I90").Value = 11700
$ws.Range("J90").Value = 32999.85
$ws.Range("K90").Value = 105300
$ws.Range("L90").Value = 296998.65
$ws.Range("M90").Value = -99060
$ws.Range("N90").Value = -309478.65

# Row 129
$ws.Range("H129").Value = 3790391
$ws.Range("J129").Value = 5558100
$ws.Range("L129").Value = 16674300
$ws.Range("N129").Value = -16684300

# Row 131
$ws.Range("H131").Value = 782.3077
$ws.Range("J131").Value = 1066.6666
$ws.Range("L131").Value = 3199.9998
$ws.Range("N131").Value = -13279.9998

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 2889
$ws.Range("I122").Value = 2132.3333
$ws.Range("K122").Value = 6396.999899999999
$ws.Range("M122").Value = -3946.999899999999

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 1524.8
$ws.Range("J100").Value = 1973.75
$ws.Range("L100").Value = 1973.75
$ws.Range("N100").Value = -3055.75

$ws = $wb.Worksheets.Item("WVR")
# Row 135
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140

